$d = $word.ActiveDocument

# --- Insertion 1: new bullet under "Database Design" section (inserted
# right after the paragraph about elevator capacity / max capacity field,
# and right before the "Future Features" Heading2 paragraph) ---
$capacityPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Since elevator capacity could be specified*") {
        $capacityPara = $p
        break
    }
}

$capacityPara.Range.InsertParagraphAfter() | Out-Null
$capacityIdx = $capacityPara.Index
$newPara1 = $d.Paragraphs.Item($capacityIdx + 1)
$newPara1.Range.Text = "Model and serial number are not specified in the requirements but it" + [char]0x2019 + "s very likely to be relevant; it is easier in the long run to implement the fields in the models now and not use them than update all the models and CRUD operations later. In a real life scenario, this business rule would be verified with the client or Product Owner during grooming prior to implementation."

# --- Insertion 2: new bullet under "Future Features" section, appended
# after the last paragraph in the document ("A logging mechanism...") ---
$loggingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "A logging mechanism that can be switched*") {
        $loggingPara = $p
        break
    }
}

$loggingPara.Range.InsertParagraphAfter() | Out-Null
$loggingIdx = $loggingPara.Index
$newPara2 = $d.Paragraphs.Item($loggingIdx + 1)
$newPara2.Range.Text = "Since there is no requirement for database management such as archiving, deleting an elevator currently does a hard delete and removes it permanently from the database. In future, they could be soft deleted instead (moved to an archiving table, marked deleted etc). "

Write-Output "done"
